$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the style of the other header cells.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the time_taken values for each data row (F2:F11) as plain text.
$times = @(
    "2021-10-05 10:52:25.001978",
    "2021-10-05 10:52:25.001987",
    "2021-10-05 10:52:25.001991",
    "2021-10-05 10:52:25.001994",
    "2021-10-05 10:52:25.001997",
    "2021-10-05 10:52:25.002000",
    "2021-10-05 10:52:25.002002",
    "2021-10-05 10:52:25.002005",
    "2021-10-05 10:52:25.002007",
    "2021-10-05 10:52:25.002010"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
